$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vegfa"
$ws.Cells.Item(2, 3).Value = "Kdr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.454188
$ws.Cells.Item(2, 8).Value = 4.362564
$ws.Cells.Item(2, 9).Value = 0.01142919054758268
$ws.Cells.Item(2, 10).Value = 0.01142919054758268
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 139.3946303333333
$ws.Cells.Item(2, 14).Value = 418.183891
$ws.Cells.Item(2, 15).Value = 0.9207771771472824
$ws.Cells.Item(2, 16).Value = 0.9207771771472822
$ws.Cells.Item(2, 17).Value = 202.7059986951693
$ws.Cells.Item(2, 18).Value = 1824.353988256524
$ws.Cells.Item(2, 19).Value = 0.01052373780948158
$ws.Cells.Item(2, 20).Value = 0.01052373780948158

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vegfa"
$ws.Cells.Item(3, 3).Value = "Kdr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.454188
$ws.Cells.Item(3, 8).Value = 4.362564
$ws.Cells.Item(3, 9).Value = 0.01142919054758268
$ws.Cells.Item(3, 10).Value = 0.01142919054758268
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.050239
$ws.Cells.Item(3, 14).Value = 0.150717
$ws.Cells.Item(3, 15).Value = 0.0003318558576616883
$ws.Cells.Item(3, 16).Value = 0.0003318558576616882
$ws.Cells.Item(3, 17).Value = 0.073056950932
$ws.Cells.Item(3, 18).Value = 0.6575125583879999
$ws.Cells.Item(3, 19).Value = 0.000003792843831546912
$ws.Cells.Item(3, 20).Value = 0.000003792843831546911

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vegfa"
$ws.Cells.Item(4, 3).Value = "Kdr"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.454188
$ws.Cells.Item(4, 8).Value = 4.362564
$ws.Cells.Item(4, 9).Value = 0.01142919054758268
$ws.Cells.Item(4, 10).Value = 0.01142919054758268
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.220699666666667
$ws.Cells.Item(4, 14).Value = 12.662099
$ws.Cells.Item(4, 15).Value = 0.02788001170035368
$ws.Cells.Item(4, 16).Value = 0.02788001170035368
$ws.Cells.Item(4, 17).Value = 6.137690806870666
$ws.Cells.Item(4, 18).Value = 55.239217261836
$ws.Cells.Item(4, 19).Value = 0.0003186459661921768
$ws.Cells.Item(4, 20).Value = 0.0003186459661921768

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vegfa"
$ws.Cells.Item(5, 3).Value = "Kdr"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.454188
$ws.Cells.Item(5, 8).Value = 4.362564
$ws.Cells.Item(5, 9).Value = 0.01142919054758268
$ws.Cells.Item(5, 10).Value = 0.01142919054758268
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.099343
$ws.Cells.Item(5, 14).Value = 18.298029
$ws.Cells.Item(5, 15).Value = 0.04028947038033828
$ws.Cells.Item(5, 16).Value = 0.04028947038033828
$ws.Cells.Item(5, 17).Value = 8.869591398484001
$ws.Cells.Item(5, 18).Value = 79.826322586356
$ws.Cells.Item(5, 19).Value = 0.0004604760340380747
$ws.Cells.Item(5, 20).Value = 0.0004604760340380746

$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Vegfa"
$ws.Cells.Item(6, 3).Value = "Kdr"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.454188
$ws.Cells.Item(6, 8).Value = 4.362564
$ws.Cells.Item(6, 9).Value = 0.01142919054758268
$ws.Cells.Item(6, 10).Value = 0.01142919054758268
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8278236666666666
$ws.Cells.Item(6, 14).Value = 2.483471
$ws.Cells.Item(6, 15).Value = 0.005468224544563193
$ws.Cells.Item(6, 16).Value = 0.005468224544563191
$ws.Cells.Item(6, 17).Value = 1.203811242182667
$ws.Cells.Item(6, 18).Value = 10.834301179644
$ws.Cells.Item(6, 19).Value = 0.00006249738027678125
$ws.Cells.Item(6, 20).Value = 0.00006249738027678124

$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Vegfa"
$ws.Cells.Item(7, 3).Value = "Kdr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.454188
$ws.Cells.Item(7, 8).Value = 4.362564
$ws.Cells.Item(7, 9).Value = 0.01142919054758268
$ws.Cells.Item(7, 10).Value = 0.01142919054758268
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7952806666666667
$ws.Cells.Item(7, 14).Value = 2.385842
$ws.Cells.Item(7, 15).Value = 0.005253260369800871
$ws.Cells.Item(7, 16).Value = 0.00525326036980087
$ws.Cells.Item(7, 17).Value = 1.156487602098667
$ws.Cells.Item(7, 18).Value = 10.408388418888
$ws.Cells.Item(7, 19).Value = 0.00006004051376251881
$ws.Cells.Item(7, 20).Value = 0.0000600405137625188

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Vegfa"
$ws.Cells.Item(8, 3).Value = "Kdr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 22.79785
$ws.Cells.Item(8, 8).Value = 68.39354999999999
$ws.Cells.Item(8, 9).Value = 0.1791797014727173
$ws.Cells.Item(8, 10).Value = 0.1791797014727173
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 139.3946303333333
$ws.Cells.Item(8, 14).Value = 418.183891
$ws.Cells.Item(8, 15).Value = 0.9207771771472824
$ws.Cells.Item(8, 16).Value = 0.9207771771472822
$ws.Cells.Item(8, 17).Value = 3177.897873144782
$ws.Cells.Item(8, 18).Value = 28601.08085830304
$ws.Cells.Item(8, 19).Value = 0.1649845797241414
$ws.Cells.Item(8, 20).Value = 0.1649845797241414

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Vegfa"
$ws.Cells.Item(9, 3).Value = "Kdr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 22.79785
$ws.Cells.Item(9, 8).Value = 68.39354999999999
$ws.Cells.Item(9, 9).Value = 0.1791797014727173
$ws.Cells.Item(9, 10).Value = 0.1791797014727173
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.050239
$ws.Cells.Item(9, 14).Value = 0.150717
$ws.Cells.Item(9, 15).Value = 0.0003318558576616883
$ws.Cells.Item(9, 16).Value = 0.0003318558576616882
$ws.Cells.Item(9, 17).Value = 1.14534118615
$ws.Cells.Item(9, 18).Value = 10.30807067535
$ws.Cells.Item(9, 19).Value = 0.00005946183350779386
$ws.Cells.Item(9, 20).Value = 0.00005946183350779386

$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Vegfa"
$ws.Cells.Item(10, 3).Value = "Kdr"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 22.79785
$ws.Cells.Item(10, 8).Value = 68.39354999999999
$ws.Cells.Item(10, 9).Value = 0.1791797014727173
$ws.Cells.Item(10, 10).Value = 0.1791797014727173
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.220699666666667
$ws.Cells.Item(10, 14).Value = 12.662099
$ws.Cells.Item(10, 15).Value = 0.02788001170035368
$ws.Cells.Item(10, 16).Value = 0.02788001170035368
$ws.Cells.Item(10, 17).Value = 96.22287789571665
$ws.Cells.Item(10, 18).Value = 866.0059010614499
$ws.Cells.Item(10, 19).Value = 0.004995532173525237
$ws.Cells.Item(10, 20).Value = 0.004995532173525237

$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Vegfa"
$ws.Cells.Item(11, 3).Value = "Kdr"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 22.79785
$ws.Cells.Item(11, 8).Value = 68.39354999999999
$ws.Cells.Item(11, 9).Value = 0.1791797014727173
$ws.Cells.Item(11, 10).Value = 0.1791797014727173
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 6.099343
$ws.Cells.Item(11, 14).Value = 18.298029
$ws.Cells.Item(11, 15).Value = 0.04028947038033828
$ws.Cells.Item(11, 16).Value = 0.04028947038033828
$ws.Cells.Item(11, 17).Value = 139.05190681255
$ws.Cells.Item(11, 18).Value = 1251.46716131295
$ws.Cells.Item(11, 19).Value = 0.007219055275242898
$ws.Cells.Item(11, 20).Value = 0.007219055275242897

$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Vegfa"
$ws.Cells.Item(12, 3).Value = "Kdr"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 22.79785
$ws.Cells.Item(12, 8).Value = 68.39354999999999
$ws.Cells.Item(12, 9).Value = 0.1791797014727173
$ws.Cells.Item(12, 10).Value = 0.1791797014727173
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.8278236666666666
$ws.Cells.Item(12, 14).Value = 2.483471
$ws.Cells.Item(12, 15).Value = 0.005468224544563193
$ws.Cells.Item(12, 16).Value = 0.005468224544563191
$ws.Cells.Item(12, 17).Value = 18.87259977911666
$ws.Cells.Item(12, 18).Value = 169.85339801205
$ws.Cells.Item(12, 19).Value = 0.0009797948414806182
$ws.Cells.Item(12, 20).Value = 0.000979794841480618

$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Vegfa"
$ws.Cells.Item(13, 3).Value = "Kdr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 22.79785
$ws.Cells.Item(13, 8).Value = 68.39354999999999
$ws.Cells.Item(13, 9).Value = 0.1791797014727173
$ws.Cells.Item(13, 10).Value = 0.1791797014727173
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.7952806666666667
$ws.Cells.Item(13, 14).Value = 2.385842
$ws.Cells.Item(13, 15).Value = 0.005253260369800871
$ws.Cells.Item(13, 16).Value = 0.00525326036980087
$ws.Cells.Item(13, 17).Value = 18.13068934656667
$ws.Cells.Item(13, 18).Value = 163.1762041191
$ws.Cells.Item(13, 19).Value = 0.0009412776248193764
$ws.Cells.Item(13, 20).Value = 0.0009412776248193762

$ws.Cells.Item(14, 1).Value = "M1"
$ws.Cells.Item(14, 2).Value = "Vegfa"
$ws.Cells.Item(14, 3).Value = "Kdr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 11.79385566666667
$ws.Cells.Item(14, 8).Value = 35.381567
$ws.Cells.Item(14, 9).Value = 0.09269380830059187
$ws.Cells.Item(14, 10).Value = 0.09269380830059187
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 139.3946303333333
$ws.Cells.Item(14, 14).Value = 418.183891
$ws.Cells.Item(14, 15).Value = 0.9207771771472824
$ws.Cells.Item(14, 16).Value = 0.9207771771472822
$ws.Cells.Item(14, 17).Value = 1644.000150859688
$ws.Cells.Item(14, 18).Value = 14796.00135773719
$ws.Cells.Item(14, 19).Value = 0.08535034314605032
$ws.Cells.Item(14, 20).Value = 0.0853503431460503

$ws.Cells.Item(15, 1).Value = "M1"
$ws.Cells.Item(15, 2).Value = "Vegfa"
$ws.Cells.Item(15, 3).Value = "Kdr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 11.79385566666667
$ws.Cells.Item(15, 8).Value = 35.381567
$ws.Cells.Item(15, 9).Value = 0.09269380830059187
$ws.Cells.Item(15, 10).Value = 0.09269380830059187
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.050239
$ws.Cells.Item(15, 14).Value = 0.150717
$ws.Cells.Item(15, 15).Value = 0.0003318558576616883
$ws.Cells.Item(15, 16).Value = 0.0003318558576616882
$ws.Cells.Item(15, 17).Value = 0.5925115148376666
$ws.Cells.Item(15, 18).Value = 5.332603633538999
$ws.Cells.Item(15, 19).Value = 0.00003076098325352104
$ws.Cells.Item(15, 20).Value = 0.00003076098325352103

$ws.Cells.Item(16, 1).Value = "M1"
$ws.Cells.Item(16, 2).Value = "Vegfa"
$ws.Cells.Item(16, 3).Value = "Kdr"
$ws.Cells.Item(16, 4).Value = "M1"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 11.79385566666667
$ws.Cells.Item(16, 8).Value = 35.381567
$ws.Cells.Item(16, 9).Value = 0.09269380830059187
$ws.Cells.Item(16, 10).Value = 0.09269380830059187
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 4.220699666666667
$ws.Cells.Item(16, 14).Value = 12.662099
$ws.Cells.Item(16, 15).Value = 0.02788001170035368
$ws.Cells.Item(16, 16).Value = 0.02788001170035368
$ws.Cells.Item(16, 17).Value = 49.77832268101477
$ws.Cells.Item(16, 18).Value = 448.004904129133
$ws.Cells.Item(16, 19).Value = 0.002584304459970843
$ws.Cells.Item(16, 20).Value = 0.002584304459970842

$ws.Cells.Item(17, 1).Value = "M1"
$ws.Cells.Item(17, 2).Value = "Vegfa"
$ws.Cells.Item(17, 3).Value = "Kdr"
$ws.Cells.Item(17, 4).Value = "M2"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 11.79385566666667
$ws.Cells.Item(17, 8).Value = 35.381567
$ws.Cells.Item(17, 9).Value = 0.09269380830059187
$ws.Cells.Item(17, 10).Value = 0.09269380830059187
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 6.099343
$ws.Cells.Item(17, 14).Value = 18.298029
$ws.Cells.Item(17, 15).Value = 0.04028947038033828
$ws.Cells.Item(17, 16).Value = 0.04028947038033828
$ws.Cells.Item(17, 17).Value = 71.93477100349367
$ws.Cells.Item(17, 18).Value = 647.4129390314429
$ws.Cells.Item(17, 19).Value = 0.003734584443967451
$ws.Cells.Item(17, 20).Value = 0.00373458444396745

$ws.Cells.Item(18, 1).Value = "M1"
$ws.Cells.Item(18, 2).Value = "Vegfa"
$ws.Cells.Item(18, 3).Value = "Kdr"
$ws.Cells.Item(18, 4).Value = "Neutro"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 11.79385566666667
$ws.Cells.Item(18, 8).Value = 35.381567
$ws.Cells.Item(18, 9).Value = 0.09269380830059187
$ws.Cells.Item(18, 10).Value = 0.09269380830059187
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.8278236666666666
$ws.Cells.Item(18, 14).Value = 2.483471
$ws.Cells.Item(18, 15).Value = 0.005468224544563193
$ws.Cells.Item(18, 16).Value = 0.005468224544563191
$ws.Cells.Item(18, 17).Value = 9.763232842117443
$ws.Cells.Item(18, 18).Value = 87.86909557905699
$ws.Cells.Item(18, 19).Value = 0.0005068705576783319
$ws.Cells.Item(18, 20).Value = 0.0005068705576783318

$ws.Cells.Item(19, 1).Value = "M1"
$ws.Cells.Item(19, 2).Value = "Vegfa"
$ws.Cells.Item(19, 3).Value = "Kdr"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 11.79385566666667
$ws.Cells.Item(19, 8).Value = 35.381567
$ws.Cells.Item(19, 9).Value = 0.09269380830059187
$ws.Cells.Item(19, 10).Value = 0.09269380830059187
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.7952806666666667
$ws.Cells.Item(19, 14).Value = 2.385842
$ws.Cells.Item(19, 15).Value = 0.005253260369800871
$ws.Cells.Item(19, 16).Value = 0.00525326036980087
$ws.Cells.Item(19, 17).Value = 9.379425397157112
$ws.Cells.Item(19, 18).Value = 84.414828574414
$ws.Cells.Item(19, 19).Value = 0.0004869447096714183
$ws.Cells.Item(19, 20).Value = 0.0004869447096714182

$ws.Cells.Item(20, 1).Value = "M2"
$ws.Cells.Item(20, 2).Value = "Vegfa"
$ws.Cells.Item(20, 3).Value = "Kdr"
$ws.Cells.Item(20, 4).Value = "ECs"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 11.85926133333333
$ws.Cells.Item(20, 8).Value = 35.577784
$ws.Cells.Item(20, 9).Value = 0.09320786413603065
$ws.Cells.Item(20, 10).Value = 0.09320786413603065
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 139.3946303333333
$ws.Cells.Item(20, 14).Value = 418.183891
$ws.Cells.Item(20, 15).Value = 0.9207771771472824
$ws.Cells.Item(20, 16).Value = 0.9207771771472822
$ws.Cells.Item(20, 17).Value = 1653.117349586394
$ws.Cells.Item(20, 18).Value = 14878.05614627754
$ws.Cells.Item(20, 19).Value = 0.08582367402710171
$ws.Cells.Item(20, 20).Value = 0.08582367402710171

$ws.Cells.Item(21, 1).Value = "M2"
$ws.Cells.Item(21, 2).Value = "Vegfa"
$ws.Cells.Item(21, 3).Value = "Kdr"
$ws.Cells.Item(21, 4).Value = "FAPs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 11.85926133333333
$ws.Cells.Item(21, 8).Value = 35.577784
$ws.Cells.Item(21, 9).Value = 0.09320786413603065
$ws.Cells.Item(21, 10).Value = 0.09320786413603065
$ws.Cells.Item(21, 11).Value = 1
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.050239
$ws.Cells.Item(21, 14).Value = 0.150717
$ws.Cells.Item(21, 15).Value = 0.0003318558576616883
$ws.Cells.Item(21, 16).Value = 0.0003318558576616882
$ws.Cells.Item(21, 17).Value = 0.5957974301253334
$ws.Cells.Item(21, 18).Value = 5.362176871128
$ws.Cells.Item(21, 19).Value = 0.00003093157569367657
$ws.Cells.Item(21, 20).Value = 0.00003093157569367656

$ws.Cells.Item(22, 1).Value = "M2"
$ws.Cells.Item(22, 2).Value = "Vegfa"
$ws.Cells.Item(22, 3).Value = "Kdr"
$ws.Cells.Item(22, 4).Value = "M1"
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 11.85926133333333
$ws.Cells.Item(22, 8).Value = 35.577784
$ws.Cells.Item(22, 9).Value = 0.09320786413603065
$ws.Cells.Item(22, 10).Value = 0.09320786413603065
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 4.220699666666667
$ws.Cells.Item(22, 14).Value = 12.662099
$ws.Cells.Item(22, 15).Value = 0.02788001170035368
$ws.Cells.Item(22, 16).Value = 0.02788001170035368
$ws.Cells.Item(22, 17).Value = 50.05438035651289
$ws.Cells.Item(22, 18).Value = 450.489423208616
$ws.Cells.Item(22, 19).Value = 0.002598636342677511
$ws.Cells.Item(22, 20).Value = 0.00259863634267751

$ws.Cells.Item(23, 1).Value = "M2"
$ws.Cells.Item(23, 2).Value = "Vegfa"
$ws.Cells.Item(23, 3).Value = "Kdr"
$ws.Cells.Item(23, 4).Value = "M2"
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 11.85926133333333
$ws.Cells.Item(23, 8).Value = 35.577784
$ws.Cells.Item(23, 9).Value = 0.09320786413603065
$ws.Cells.Item(23, 10).Value = 0.09320786413603065
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 6.099343
$ws.Cells.Item(23, 14).Value = 18.298029
$ws.Cells.Item(23, 15).Value = 0.04028947038033828
$ws.Cells.Item(23, 16).Value = 0.04028947038033828
$ws.Cells.Item(23, 17).Value = 72.33370259863734
$ws.Cells.Item(23, 18).Value = 651.003323387736
$ws.Cells.Item(23, 19).Value = 0.003755295481323201
$ws.Cells.Item(23, 20).Value = 0.003755295481323201

$ws.Cells.Item(24, 1).Value = "M2"
$ws.Cells.Item(24, 2).Value = "Vegfa"
$ws.Cells.Item(24, 3).Value = "Kdr"
$ws.Cells.Item(24, 4).Value = "Neutro"
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 11.85926133333333
$ws.Cells.Item(24, 8).Value = 35.577784
$ws.Cells.Item(24, 9).Value = 0.09320786413603065
$ws.Cells.Item(24, 10).Value = 0.09320786413603065
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 0.8278236666666666
$ws.Cells.Item(24, 14).Value = 2.483471
$ws.Cells.Item(24, 15).Value = 0.005468224544563193
$ws.Cells.Item(24, 16).Value = 0.005468224544563191
$ws.Cells.Item(24, 17).Value = 9.817377200918223
$ws.Cells.Item(24, 18).Value = 88.356394808264
$ws.Cells.Item(24, 19).Value = 0.0005096815304149541
$ws.Cells.Item(24, 20).Value = 0.000509681530414954

$ws.Cells.Item(25, 1).Value = "M2"
$ws.Cells.Item(25, 2).Value = "Vegfa"
$ws.Cells.Item(25, 3).Value = "Kdr"
$ws.Cells.Item(25, 4).Value = "sCs"
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 11.85926133333333
$ws.Cells.Item(25, 8).Value = 35.577784
$ws.Cells.Item(25, 9).Value = 0.09320786413603065
$ws.Cells.Item(25, 10).Value = 0.09320786413603065
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 0.7952806666666667
$ws.Cells.Item(25, 14).Value = 2.385842
$ws.Cells.Item(25, 15).Value = 0.005253260369800871
$ws.Cells.Item(25, 16).Value = 0.00525326036980087
$ws.Cells.Item(25, 17).Value = 9.431441259347558
$ws.Cells.Item(25, 18).Value = 84.88297133412802
$ws.Cells.Item(25, 19).Value = 0.0004896451788195936
$ws.Cells.Item(25, 20).Value = 0.0004896451788195936

$ws.Cells.Item(26, 1).Value = "Neutro"
$ws.Cells.Item(26, 2).Value = "Vegfa"
$ws.Cells.Item(26, 3).Value = "Kdr"
$ws.Cells.Item(26, 4).Value = "ECs"
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 71.30949033333333
$ws.Cells.Item(26, 8).Value = 213.928471
$ws.Cells.Item(26, 9).Value = 0.5604569373909507
$ws.Cells.Item(26, 10).Value = 0.5604569373909508
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = 139.3946303333333
$ws.Cells.Item(26, 14).Value = 418.183891
$ws.Cells.Item(26, 15).Value = 0.9207771771472824
$ws.Cells.Item(26, 16).Value = 0.9207771771472822
$ws.Cells.Item(26, 17).Value = 9940.160044273405
$ws.Cells.Item(26, 18).Value = 89461.44039846065
$ws.Cells.Item(26, 19).Value = 0.5160559567234507
$ws.Cells.Item(26, 20).Value = 0.5160559567234508

$ws.Cells.Item(27, 1).Value = "Neutro"
$ws.Cells.Item(27, 2).Value = "Vegfa"
$ws.Cells.Item(27, 3).Value = "Kdr"
$ws.Cells.Item(27, 4).Value = "FAPs"
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 71.30949033333333
$ws.Cells.Item(27, 8).Value = 213.928471
$ws.Cells.Item(27, 9).Value = 0.5604569373909507
$ws.Cells.Item(27, 10).Value = 0.5604569373909508
$ws.Cells.Item(27, 11).Value = 1
$ws.Cells.Item(27, 12).Value = 0.3333333333333333
$ws.Cells.Item(27, 13).Value = 0.050239
$ws.Cells.Item(27, 14).Value = 0.150717
$ws.Cells.Item(27, 15).Value = 0.0003318558576616883
$ws.Cells.Item(27, 16).Value = 0.0003318558576616882
$ws.Cells.Item(27, 17).Value = 3.582517484856333
$ws.Cells.Item(27, 18).Value = 32.242657363707
$ws.Cells.Item(27, 19).Value = 0.0001859909176403171
$ws.Cells.Item(27, 20).Value = 0.0001859909176403171

$ws.Cells.Item(28, 1).Value = "Neutro"
$ws.Cells.Item(28, 2).Value = "Vegfa"
$ws.Cells.Item(28, 3).Value = "Kdr"
$ws.Cells.Item(28, 4).Value = "M1"
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 71.30949033333333
$ws.Cells.Item(28, 8).Value = 213.928471
$ws.Cells.Item(28, 9).Value = 0.5604569373909507
$ws.Cells.Item(28, 10).Value = 0.5604569373909508
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 12).Value = 1
$ws.Cells.Item(28, 13).Value = 4.220699666666667
$ws.Cells.Item(28, 14).Value = 12.662099
$ws.Cells.Item(28, 15).Value = 0.02788001170035368
$ws.Cells.Item(28, 16).Value = 0.02788001170035368
$ws.Cells.Item(28, 17).Value = 300.9759420800698
$ws.Cells.Item(28, 18).Value = 2708.783478720629
$ws.Cells.Item(28, 19).Value = 0.0156255459720041
$ws.Cells.Item(28, 20).Value = 0.0156255459720041

$ws.Cells.Item(29, 1).Value = "Neutro"
$ws.Cells.Item(29, 2).Value = "Vegfa"
$ws.Cells.Item(29, 3).Value = "Kdr"
$ws.Cells.Item(29, 4).Value = "M2"
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 71.30949033333333
$ws.Cells.Item(29, 8).Value = 213.928471
$ws.Cells.Item(29, 9).Value = 0.5604569373909507
$ws.Cells.Item(29, 10).Value = 0.5604569373909508
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 12).Value = 1
$ws.Cells.Item(29, 13).Value = 6.099343
$ws.Cells.Item(29, 14).Value = 18.298029
$ws.Cells.Item(29, 15).Value = 0.04028947038033828
$ws.Cells.Item(29, 16).Value = 0.04028947038033828
$ws.Cells.Item(29, 17).Value = 434.9410406981843
$ws.Cells.Item(29, 18).Value = 3914.469366283659
$ws.Cells.Item(29, 19).Value = 0.02258051317846782
$ws.Cells.Item(29, 20).Value = 0.02258051317846782

$ws.Cells.Item(30, 1).Value = "Neutro"
$ws.Cells.Item(30, 2).Value = "Vegfa"
$ws.Cells.Item(30, 3).Value = "Kdr"
$ws.Cells.Item(30, 4).Value = "Neutro"
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 71.30949033333333
$ws.Cells.Item(30, 8).Value = 213.928471
$ws.Cells.Item(30, 9).Value = 0.5604569373909507
$ws.Cells.Item(30, 10).Value = 0.5604569373909508
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 12).Value = 1
$ws.Cells.Item(30, 13).Value = 0.8278236666666666
$ws.Cells.Item(30, 14).Value = 2.483471
$ws.Cells.Item(30, 15).Value = 0.005468224544563193
$ws.Cells.Item(30, 16).Value = 0.005468224544563191
$ws.Cells.Item(30, 17).Value = 59.03168375587121
$ws.Cells.Item(30, 18).Value = 531.2851538028409
$ws.Cells.Item(30, 19).Value = 0.003064704381211913
$ws.Cells.Item(30, 20).Value = 0.003064704381211913

$ws.Cells.Item(31, 1).Value = "Neutro"
$ws.Cells.Item(31, 2).Value = "Vegfa"
$ws.Cells.Item(31, 3).Value = "Kdr"
$ws.Cells.Item(31, 4).Value = "sCs"
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 71.30949033333333
$ws.Cells.Item(31, 8).Value = 213.928471
$ws.Cells.Item(31, 9).Value = 0.5604569373909507
$ws.Cells.Item(31, 10).Value = 0.5604569373909508
$ws.Cells.Item(31, 11).Value = 3
$ws.Cells.Item(31, 12).Value = 1
$ws.Cells.Item(31, 13).Value = 0.7952806666666667
$ws.Cells.Item(31, 14).Value = 2.385842
$ws.Cells.Item(31, 15).Value = 0.005253260369800871
$ws.Cells.Item(31, 16).Value = 0.00525326036980087
$ws.Cells.Item(31, 17).Value = 56.71105901195356
$ws.Cells.Item(31, 18).Value = 510.399531107582
$ws.Cells.Item(31, 19).Value = 0.002944226218175849
$ws.Cells.Item(31, 20).Value = 0.002944226218175849

$ws.Cells.Item(32, 1).Value = "sCs"
$ws.Cells.Item(32, 2).Value = "Vegfa"
$ws.Cells.Item(32, 3).Value = "Kdr"
$ws.Cells.Item(32, 4).Value = "ECs"
$ws.Cells.Item(32, 5).Value = 3
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 8.019912
$ws.Cells.Item(32, 8).Value = 24.059736
$ws.Cells.Item(32, 9).Value = 0.06303249815212676
$ws.Cells.Item(32, 10).Value = 0.06303249815212676
$ws.Cells.Item(32, 11).Value = 3
$ws.Cells.Item(32, 12).Value = 1
$ws.Cells.Item(32, 13).Value = 139.3946303333333
$ws.Cells.Item(32, 14).Value = 418.183891
$ws.Cells.Item(32, 15).Value = 0.9207771771472824
$ws.Cells.Item(32, 16).Value = 0.9207771771472822
$ws.Cells.Item(32, 17).Value = 1117.932668545864
$ws.Cells.Item(32, 18).Value = 10061.39401691277
$ws.Cells.Item(32, 19).Value = 0.05803888571705657
$ws.Cells.Item(32, 20).Value = 0.05803888571705657

$ws.Cells.Item(33, 1).Value = "sCs"
$ws.Cells.Item(33, 2).Value = "Vegfa"
$ws.Cells.Item(33, 3).Value = "Kdr"
$ws.Cells.Item(33, 4).Value = "FAPs"
$ws.Cells.Item(33, 5).Value = 3
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 8.019912
$ws.Cells.Item(33, 8).Value = 24.059736
$ws.Cells.Item(33, 9).Value = 0.06303249815212676
$ws.Cells.Item(33, 10).Value = 0.06303249815212676
$ws.Cells.Item(33, 11).Value = 1
$ws.Cells.Item(33, 12).Value = 0.3333333333333333
$ws.Cells.Item(33, 13).Value = 0.050239
$ws.Cells.Item(33, 14).Value = 0.150717
$ws.Cells.Item(33, 15).Value = 0.0003318558576616883
$ws.Cells.Item(33, 16).Value = 0.0003318558576616882
$ws.Cells.Item(33, 17).Value = 0.402912358968
$ws.Cells.Item(33, 18).Value = 3.626211230712
$ws.Cells.Item(33, 19).Value = 0.00002091770373483281
$ws.Cells.Item(33, 20).Value = 0.0000209177037348328

$ws.Cells.Item(34, 1).Value = "sCs"
$ws.Cells.Item(34, 2).Value = "Vegfa"
$ws.Cells.Item(34, 3).Value = "Kdr"
$ws.Cells.Item(34, 4).Value = "M1"
$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 8.019912
$ws.Cells.Item(34, 8).Value = 24.059736
$ws.Cells.Item(34, 9).Value = 0.06303249815212676
$ws.Cells.Item(34, 10).Value = 0.06303249815212676
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 12).Value = 1
$ws.Cells.Item(34, 13).Value = 4.220699666666667
$ws.Cells.Item(34, 14).Value = 12.662099
$ws.Cells.Item(34, 15).Value = 0.02788001170035368
$ws.Cells.Item(34, 16).Value = 0.02788001170035368
$ws.Cells.Item(34, 17).Value = 33.849639905096
$ws.Cells.Item(34, 18).Value = 304.646759145864
$ws.Cells.Item(34, 19).Value = 0.001757346785983816
$ws.Cells.Item(34, 20).Value = 0.001757346785983816

$ws.Cells.Item(35, 1).Value = "sCs"
$ws.Cells.Item(35, 2).Value = "Vegfa"
$ws.Cells.Item(35, 3).Value = "Kdr"
$ws.Cells.Item(35, 4).Value = "M2"
$ws.Cells.Item(35, 5).Value = 3
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 8.019912
$ws.Cells.Item(35, 8).Value = 24.059736
$ws.Cells.Item(35, 9).Value = 0.06303249815212676
$ws.Cells.Item(35, 10).Value = 0.06303249815212676
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 12).Value = 1
$ws.Cells.Item(35, 13).Value = 6.099343
$ws.Cells.Item(35, 14).Value = 18.298029
$ws.Cells.Item(35, 15).Value = 0.04028947038033828
$ws.Cells.Item(35, 16).Value = 0.04028947038033828
$ws.Cells.Item(35, 17).Value = 48.916194117816
$ws.Cells.Item(35, 18).Value = 440.245747060344
$ws.Cells.Item(35, 19).Value = 0.002539545967298839
$ws.Cells.Item(35, 20).Value = 0.002539545967298838

$ws.Cells.Item(36, 1).Value = "sCs"
$ws.Cells.Item(36, 2).Value = "Vegfa"
$ws.Cells.Item(36, 3).Value = "Kdr"
$ws.Cells.Item(36, 4).Value = "Neutro"
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 8.019912
$ws.Cells.Item(36, 8).Value = 24.059736
$ws.Cells.Item(36, 9).Value = 0.06303249815212676
$ws.Cells.Item(36, 10).Value = 0.06303249815212676
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 12).Value = 1
$ws.Cells.Item(36, 13).Value = 0.8278236666666666
$ws.Cells.Item(36, 14).Value = 2.483471
$ws.Cells.Item(36, 15).Value = 0.005468224544563193
$ws.Cells.Item(36, 16).Value = 0.005468224544563191
$ws.Cells.Item(36, 17).Value = 6.639072958183999
$ws.Cells.Item(36, 18).Value = 59.751656623656
$ws.Cells.Item(36, 19).Value = 0.0003446758535005937
$ws.Cells.Item(36, 20).Value = 0.0003446758535005936

$ws.Cells.Item(37, 1).Value = "sCs"
$ws.Cells.Item(37, 2).Value = "Vegfa"
$ws.Cells.Item(37, 3).Value = "Kdr"
$ws.Cells.Item(37, 4).Value = "sCs"
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 8.019912
$ws.Cells.Item(37, 8).Value = 24.059736
$ws.Cells.Item(37, 9).Value = 0.06303249815212676
$ws.Cells.Item(37, 10).Value = 0.06303249815212676
$ws.Cells.Item(37, 11).Value = 3
$ws.Cells.Item(37, 12).Value = 1
$ws.Cells.Item(37, 13).Value = 0.7952806666666667
$ws.Cells.Item(37, 14).Value = 2.385842
$ws.Cells.Item(37, 15).Value = 0.005253260369800871
$ws.Cells.Item(37, 16).Value = 0.00525326036980087
$ws.Cells.Item(37, 17).Value = 6.378080961968
$ws.Cells.Item(37, 18).Value = 57.40272865771201
$ws.Cells.Item(37, 19).Value = 0.0003311261245521141
$ws.Cells.Item(37, 20).Value = 0.0003311261245521141

